# Auto-generated Excel COM-interop script
# Applies the "Add data for 2024-08-30" crime-count data refresh
# across the Citywide Totals, By Neighborhood, and per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

# --- Citywide Totals (32 cells) ---
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("C2").Value = 45
$ws.Range("K2").Value = 96
$ws.Range("H3").Value = 89
$ws.Range("I3").Value = 134
$ws.Range("J3").Value = 140
$ws.Range("C9").Value = 323
$ws.Range("J9").Value = 274
$ws.Range("H9").Value = 293
$ws.Range("G9").Value = 330
$ws.Range("D9").Value = 291
$ws.Range("E9").Value = 286
$ws.Range("F9").Value = 376
$ws.Range("I9").Value = 363
$ws.Range("G10").Value = 732
$ws.Range("I10").Value = 583
$ws.Range("C10").Value = 1003
$ws.Range("J10").Value = 477
$ws.Range("B10").Value = 842
$ws.Range("K10").Value = 495
$ws.Range("D10").Value = 1183
$ws.Range("E10").Value = 1465
$ws.Range("F10").Value = 1503
$ws.Range("B11").Value = 1188
$ws.Range("D11").Value = 1638
$ws.Range("E11").Value = 1902
$ws.Range("F11").Value = 2038
$ws.Range("H11").Value = 828
$ws.Range("G11").Value = 1215
$ws.Range("I11").Value = 1184
$ws.Range("C11").Value = 1440
$ws.Range("J11").Value = 994
$ws.Range("K11").Value = 1104

# --- Chicago Lawn (2 cells) ---
$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J5").Value = 2
$ws.Range("J7").Value = 7

# --- Garfield Park (10 cells) ---
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 4
$ws.Range("C7").Value = 20
$ws.Range("E7").Value = 27
$ws.Range("D7").Value = 29
$ws.Range("I8").Value = 37
$ws.Range("C9").Value = 58
$ws.Range("D9").Value = 70
$ws.Range("E9").Value = 91
$ws.Range("K9").Value = 54
$ws.Range("I9").Value = 72

# --- Chatham (4 cells) ---
$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I6").Value = 8
$ws.Range("G6").Value = 11
$ws.Range("G8").Value = 29
$ws.Range("I8").Value = 23

# --- Grand Crossing (4 cells) ---
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("H7").Value = 13
$ws.Range("K8").Value = 25
$ws.Range("H9").Value = 40
$ws.Range("K9").Value = 72

# --- Loop (12 cells) ---
$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I3").Value = 21
$ws.Range("C9").Value = 190
$ws.Range("D9").Value = 346
$ws.Range("E9").Value = 425
$ws.Range("B9").Value = 109
$ws.Range("F9").Value = 407
$ws.Range("B10").Value = 142
$ws.Range("D10").Value = 403
$ws.Range("E10").Value = 477
$ws.Range("F10").Value = 462
$ws.Range("C10").Value = 228
$ws.Range("I10").Value = 232

# --- Old Town (6 cells) ---
$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("H3").Value = 2
$ws.Range("E6").Value = 36
$ws.Range("G6").Value = 19
$ws.Range("E7").Value = 45
$ws.Range("H7").Value = 17
$ws.Range("G7").Value = 31

# --- Little Italy, UIC (4 cells) ---
$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("F5").Value = 21
$ws.Range("D5").Value = 6
$ws.Range("F7").Value = 50
$ws.Range("D7").Value = 21

# --- North Lawndale (4 cells) ---
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("F7").Value = 27
$ws.Range("E7").Value = 23
$ws.Range("F8").Value = 51
$ws.Range("E8").Value = 30

# --- By Neighborhood (66 cells) ---
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I6").Value = 3
$ws.Range("J10").Value = 13
$ws.Range("J11").Value = 5
$ws.Range("E16").Value = 5
$ws.Range("I19").Value = 23
$ws.Range("G19").Value = 29
$ws.Range("J20").Value = 7
$ws.Range("D22").Value = 11
$ws.Range("J22").Value = 3
$ws.Range("K27").Value = 12
$ws.Range("G28").Value = 58
$ws.Range("H29").Value = 9
$ws.Range("C30").Value = 20
$ws.Range("F31").Value = 2
$ws.Range("D32").Value = 70
$ws.Range("E32").Value = 91
$ws.Range("C32").Value = 58
$ws.Range("K32").Value = 54
$ws.Range("I32").Value = 72
$ws.Range("H36").Value = 40
$ws.Range("K36").Value = 72
$ws.Range("C43").Value = 9
$ws.Range("G48").Value = 6
$ws.Range("D50").Value = 21
$ws.Range("F50").Value = 50
$ws.Range("I51").Value = 6
$ws.Range("C52").Value = 20
$ws.Range("C53").Value = 228
$ws.Range("D53").Value = 403
$ws.Range("B53").Value = 142
$ws.Range("E53").Value = 477
$ws.Range("F53").Value = 462
$ws.Range("I53").Value = 232
$ws.Range("D57").Value = 2
$ws.Range("E61").Value = 41
$ws.Range("C62").Value = 16
$ws.Range("F65").Value = 51
$ws.Range("E65").Value = 30
$ws.Range("C66").Value = 3
$ws.Range("H70").Value = 17
$ws.Range("G70").Value = 31
$ws.Range("E70").Value = 45
$ws.Range("D72").Value = 6
$ws.Range("D74").Value = 58
$ws.Range("K74").Value = 24
$ws.Range("C76").Value = 52
$ws.Range("G76").Value = 31
$ws.Range("F76").Value = 42
$ws.Range("K77").Value = 47
$ws.Range("I86").Value = 7
$ws.Range("J87").Value = 19
$ws.Range("I87").Value = 18
$ws.Range("D89").Value = 16
$ws.Range("D92").Value = 31
$ws.Range("F96").Value = 16
$ws.Range("J96").Value = 5
$ws.Range("I99").Value = 1184
$ws.Range("C99").Value = 1440
$ws.Range("J99").Value = 994
$ws.Range("K99").Value = 1104
$ws.Range("B99").Value = 1188
$ws.Range("D99").Value = 1638
$ws.Range("E99").Value = 1902
$ws.Range("F99").Value = 2038
$ws.Range("H99").Value = 828
$ws.Range("G99").Value = 1215

# --- Washington Park (2 cells) ---
$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("D5").Value = 7
$ws.Range("D6").Value = 16

# --- Uptown (4 cells) ---
$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J8").Value = 9
$ws.Range("I8").Value = 10
$ws.Range("J9").Value = 19
$ws.Range("I9").Value = 18

# --- Englewood (2 cells) ---
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("G7").Value = 21
$ws.Range("G9").Value = 58

# --- Fuller Park (2 cells) ---
$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("H7").Value = 3
$ws.Range("H9").Value = 9

# --- Rogers Park (6 cells) ---
$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("C7").Value = 40
$ws.Range("G7").Value = 21
$ws.Range("F7").Value = 27
$ws.Range("C8").Value = 52
$ws.Range("G8").Value = 31
$ws.Range("F8").Value = 42

# --- River North (4 cells) ---
$ws = $wb.Worksheets.Item('River North')
$ws.Range("D6").Value = 48
$ws.Range("K6").Value = 13
$ws.Range("D7").Value = 58
$ws.Range("K7").Value = 24

# --- United Center (2 cells) ---
$ws = $wb.Worksheets.Item('United Center')
$ws.Range("I6").Value = 4
$ws.Range("I7").Value = 7

# --- West Loop (2 cells) ---
$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("D8").Value = 24
$ws.Range("D9").Value = 31

# --- Little Village (2 cells) ---
$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I6").Value = 5
$ws.Range("I7").Value = 6

# --- Logan Square (2 cells) ---
$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("C7").Value = 14
$ws.Range("C8").Value = 20

# --- Edgewater (2 cells) ---
$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K5").Value = 6
$ws.Range("K6").Value = 12

# --- Near South Side (2 cells) ---
$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("C6").Value = 2
$ws.Range("C8").Value = 16

# --- Millenium Park (2 cells) ---
$ws = $wb.Worksheets.Item('Millenium Park')
$ws.Range("D3").Value = 2
$ws.Range("D4").Value = 2

# --- Roseland (2 cells) ---
$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K9").Value = 27
$ws.Range("K10").Value = 47

# --- Clearing (4 cells) ---
$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("D6").Value = 9
$ws.Range("J6").Value = 2
$ws.Range("D7").Value = 11
$ws.Range("J7").Value = 3

# --- Wicker Park (4 cells) ---
$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("F5").Value = 14
$ws.Range("J5").Value = 3
$ws.Range("F6").Value = 16
$ws.Range("J6").Value = 5

# --- Lincoln Park (2 cells) ---
$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("G6").Value = 5
$ws.Range("G7").Value = 6

# --- Printers Row (3 cells) ---
$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range("D5").Value = 3
$ws.Range("D6").Value = 6
$ws.Range("J3").Value = 2

# --- Avondale (1 cell) ---
$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("J7").Value = 13

# --- Gage Park (2 cells) ---
$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("C7").Value = 17
$ws.Range("C8").Value = 20

# --- Belmont Cragin (2 cells) ---
$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I5").Value = 3
$ws.Range("I7").Value = 5

# --- Bucktown (2 cells) ---
$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("D5").Value = 4
$ws.Range("D6").Value = 5

# --- Ashburn (2 cells) ---
$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("I5").Value = 1
$ws.Range("I6").Value = 3

# --- Irving Park (2 cells) ---
$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("C6").Value = 5
$ws.Range("C7").Value = 9

# --- North Park (2 cells) ---
$ws = $wb.Worksheets.Item('North Park')
$ws.Range("C2").Value = 1
$ws.Range("C6").Value = 3

# --- Galewood (2 cells) ---
$ws = $wb.Worksheets.Item('Galewood')
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2

Write-Output "Updated 208 cells across 36 sheets"
